# Update handback status report timestamps (regenerated report values).
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for bffc659a...md (row 2, col G)
$wsOverview.Range("G2").Value = "2016-08-25 21:06:26"

# zh-cn sheet: Correspond Handoff Datetime (H2) / Correspond Handback DateTime (K2) for bffc659a row
$wsZhCn.Range("H2").Value = "2016-08-25 21:06:22"
$wsZhCn.Range("K2").Value = "2016-08-25 21:06:39"

# de-de sheet: Correspond Handback DateTime (K2) for bffc659a row
$wsDeDe.Range("K2").Value = "2016-08-25 21:06:46"
